$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.024.10'
$ws.Range('E2').Value = '  +4.68%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.416.09'
$ws.Range('E3').Value = '  +3.66%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('E5').Value = '  +4.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.56'
$ws.Range('E6').Value = '  +8.91%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.414.96'
$ws.Range('E8').Value = '  +3.66%  '
$ws.Range('E9').Value = '  +2.44%  '
$ws.Range('E10').Value = '  +2.45%  '
$ws.Range('E11').Value = '  +10.12%  '
$ws.Range('E12').Value = '  +6.91%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.003.10'
$ws.Range('E13').Value = '  +3.83%  '
$ws.Range('E14').Value = '  +2.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000179'
$ws.Range('E15').Value = '  +8.52%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.415.03'
$ws.Range('E16').Value = '  +3.78%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.47'
$ws.Range('E17').Value = '  +6.30%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '62.021.11'
$ws.Range('E18').Value = '  +4.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.16'
$ws.Range('E19').Value = '  +7.66%  '
$ws.Range('E20').Value = '  +5.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.53'
$ws.Range('E21').Value = '  +7.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '391.10'
$ws.Range('E22').Value = '  +11.96%  '
$ws.Range('E23').Value = '  +3.96%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.553.58'
$ws.Range('E24').Value = '  +3.67%  '
$ws.Range('E25').Value = '  +17.54%  '
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '71.54'
$ws.Range('E27').Value = '  +4.85%  '
$ws.Range('E28').Value = '  +6.83%  '
$ws.Range('E29').Value = '  +9.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.08%  '
$ws.Range('E31').Value = '  +7.36%  '
$ws.Range('E32').Value = '  +5.87%  '
$ws.Range('E33').Value = '  +4.87%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.446.84'
$ws.Range('E34').Value = '  +3.59%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.56'
$ws.Range('E36').Value = '  +4.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.48'
$ws.Range('E37').Value = '  +3.38%  '
$ws.Range('E38').Value = '  +4.01%  '
$ws.Range('E39').Value = '  +6.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '163.14'
$ws.Range('E40').Value = '  +3.61%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0792'
$ws.Range('E41').Value = '  +6.84%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.77'
$ws.Range('E42').Value = '  +15.99%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.790'
$ws.Range('E43').Value = '  +7.31%  '
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('E45').Value = '  +5.83%  '
$ws.Range('E46').Value = '  +5.13%  '
$ws.Range('B47').Value = 'OKB'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '41.78'
$ws.Range('E47').Value = '  +2.93%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '25.07'
$ws.Range('E48').Value = '  +11.20%  '
$ws.Range('E50').Value = '  +7.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.380.53'
$ws.Range('E51').Value = '  +11.57%  '
